$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents() | Out-Null
$ws.Range("H43").Value = 1090.4
$ws.Range("I43").Value = 983.3333
$ws.Range("K43").Value = 983.3333
$ws.Range("M43").Value = -914.3333
$ws.Range("K53").Value = 1337.1666
$ws.Range("I53").Value = 1337.1666
$ws.Range("H53").Value = 1349.3077
$ws.Range("M53").Value = -700.1666
$ws.Range("M80").Value = -839.5
$ws.Range("I80").Value = 612.5
$ws.Range("H80").Value = 966.5
$ws.Range("K80").Value = 1837.5
$ws.Range("I83").Value = 612.5
$ws.Range("M83").Value = -520.5
$ws.Range("K83").Value = 5512.5
$ws.Range("H83").Value = 966.5
$ws.Range("I132").Value = 9794.272000000001
$ws.Range("L132").Value = 19375.9995
$ws.Range("M132").Value = -26852.816
$ws.Range("K132").Value = 29382.816
$ws.Range("H132").Value = 9079.5
$ws.Range("N132").Value = -24435.9995
$ws.Range("J132").Value = 6458.6665
$ws.Range("H138").Value = 2425
$ws.Range("K138").Value = 6580.5
$ws.Range("I138").Value = 2193.5
$ws.Range("M138").Value = -1440.5
$ws.Range("J139").Value = 100780
$ws.Range("N139").Value = -111060
$ws.Range("H139").Value = 100780
$ws.Range("L139").Value = 100780

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I6").Value = 112500000
$ws.Range("M6").Value = -112499827
$ws.Range("H6").Value = 112500000
$ws.Range("K6").Value = 112500000
$ws.Range("I32").Value = 5853.6665
$ws.Range("M32").Value = -5566.6665
$ws.Range("K32").Value = 5853.6665
$ws.Range("H32").Value = 5853.6665
$ws.Range("H92").Value = 57500
$ws.Range("L92").Value = 57500
$ws.Range("J92").Value = 57500
$ws.Range("N92").Value = -62492
$ws.Range("K97").Value = 335
$ws.Range("I97").Value = 335
$ws.Range("H97").Value = 406.5
$ws.Range("M97").Value = 161
$ws.Range("H122").Value = 1580.1666
$ws.Range("K122").Value = 4740.4998
$ws.Range("I122").Value = 1580.1666
$ws.Range("M122").Value = -2290.4998
$ws.Range("I132").Value = 2832.375
$ws.Range("M132").Value = -5967.125
$ws.Range("K132").Value = 8497.125
$ws.Range("H132").Value = 4254.9165

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents() | Out-Null
$ws.Range("H46").Value = 0
$ws.Range("H64").Value = 1277
$ws.Range("I64").Value = 0
$ws.Range("L64").Value = 1277
$ws.Range("K64").Value = 0
$ws.Range("J64").Value = 1277
$ws.Range("M64").ClearContents() | Out-Null
$ws.Range("N64").Value = -1727
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 1277
$ws.Range("M67").ClearContents() | Out-Null
$ws.Range("N67").Value = -2837
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1277
$ws.Range("H67").Value = 1277
$ws.Range("H94").Value = 1222.5
$ws.Range("M94").Value = -796.2221999999999
$ws.Range("K94").Value = 1247.2222
$ws.Range("I94").Value = 1247.2222
$ws.Range("I134").Value = 2172.7144
$ws.Range("H134").Value = 2172.7144
$ws.Range("K134").Value = 6518.1432
$ws.Range("M134").Value = -3983.1432

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6189.4443
$ws.Range("M62").Value = -5717
$ws.Range("N62").Value = -7248
$ws.Range("J62").Value = 6000
$ws.Range("I62").Value = 6341
$ws.Range("K62").Value = 6341
$ws.Range("L62").Value = 6000
$ws.Range("I65").Value = 6341
$ws.Range("H65").Value = 6189.4443
$ws.Range("N65").Value = -36240
$ws.Range("M65").Value = -28585
$ws.Range("L65").Value = 30000
$ws.Range("K65").Value = 31705
$ws.Range("J65").Value = 6000

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L140").Value = 15000
$ws.Range("N140").Value = -25360
$ws.Range("H140").Value = 3666.3333
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 8998.5
$ws.Range("M140").Value = -3818.5
$ws.Range("I140").Value = 2999.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K97").Value = 676.4
$ws.Range("I97").Value = 676.4
$ws.Range("H97").Value = 751.2727
$ws.Range("M97").Value = -180.4
$ws.Range("H122").Value = 1564.6364
$ws.Range("K122").Value = 5070.3333
$ws.Range("I122").Value = 1690.1111
$ws.Range("M122").Value = -2620.3333
$ws.Range("I126").Value = 0
$ws.Range("M126").ClearContents() | Out-Null
$ws.Range("K126").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I132").Value = 835.625
$ws.Range("L132").Value = 2100
$ws.Range("M132").Value = 23.125
$ws.Range("K132").Value = 2506.875
$ws.Range("H132").Value = 820.55554
$ws.Range("N132").Value = -7160
$ws.Range("J132").Value = 700

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L7").Value = 0
$ws.Range("H7").Value = 2849.5
$ws.Range("N7").ClearContents() | Out-Null
$ws.Range("J7").Value = 0
$ws.Range("I22").Value = 530.8333
$ws.Range("K22").Value = 530.8333
$ws.Range("J22").Value = 670.8
$ws.Range("H22").Value = 594.4545000000001
$ws.Range("N22").Value = -1260.8
$ws.Range("L22").Value = 670.8
$ws.Range("M22").Value = -235.8333
$ws.Range("J27").Value = 670.8
$ws.Range("M27").Value = -423.8333
$ws.Range("N27").Value = -884.8
$ws.Range("H27").Value = 594.4545000000001
$ws.Range("L27").Value = 670.8
$ws.Range("K27").Value = 530.8333
$ws.Range("I27").Value = 530.8333
$ws.Range("H43").Value = 23750
$ws.Range("I43").Value = 23750
$ws.Range("K43").Value = 23750
$ws.Range("M43").Value = -23557
$ws.Range("H82").Value = 2022
$ws.Range("M82").Value = -1710.75
$ws.Range("I82").Value = 2071.75
$ws.Range("K82").Value = 2071.75
$ws.Range("M85").Value = -823.75
$ws.Range("H85").Value = 2022
$ws.Range("I85").Value = 2071.75
$ws.Range("K85").Value = 2071.75
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
$ws.Range("H125").Value = 80000
$ws.Range("L126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("N126").ClearContents() | Out-Null
$ws.Range("H126").Value = 2849.5
$ws.Range("I132").Value = 1919
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3227
$ws.Range("K132").Value = 5757
$ws.Range("H132").Value = 2314.25
$ws.Range("N132").Value = -15560
$ws.Range("J132").Value = 3500
$ws.Range("H136").Value = 1698
$ws.Range("M136").Value = -2544
$ws.Range("I136").Value = 1698
$ws.Range("K136").Value = 5094

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents() | Out-Null
$ws.Range("J32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H100").Value = 330.66666
$ws.Range("I100").Value = 330.66666
$ws.Range("M100").Value = -120.33332
$ws.Range("K100").Value = 661.33332
$ws.Range("I126").Value = 5522.625
$ws.Range("M126").Value = -14097.875
$ws.Range("K126").Value = 16567.875
$ws.Range("H126").Value = 5522.625
$ws.Range("H136").Value = 8633.5
$ws.Range("M136").Value = -23350.5
$ws.Range("I136").Value = 8633.5
$ws.Range("K136").Value = 25900.5
